$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast-error statistics (naive component forecaster bug fix)
# Columns: B=ME, C=MAE, D=MSE, E=RMSE, F=SE, G=N

$ws.Range("B2").Value = 0.006971674199201135
$ws.Range("C2").Value = 0.414613479980299
$ws.Range("D2").Value = 0.4465464471237033
$ws.Range("E2").Value = 0.6682413090521292
$ws.Range("F2").Value = 0.6748539099243737
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.4355859266737307
$ws.Range("C3").Value = 0.699453688851532
$ws.Range("D3").Value = 1.429667568063306
$ws.Range("E3").Value = 1.195687069455594
$ws.Range("F3").Value = 1.124827636930724
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.5940129335889834
$ws.Range("C4").Value = 1.050797193125951
$ws.Range("D4").Value = 4.135615467765216
$ws.Range("E4").Value = 2.033621269500596
$ws.Range("F4").Value = 1.965088213786561
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.5024337226471743
$ws.Range("C5").Value = 1.099887659833765
$ws.Range("D5").Value = 4.742322585733941
$ws.Range("E5").Value = 2.177687439862282
$ws.Range("F5").Value = 2.141357598426261
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.4027908244544923
$ws.Range("C6").Value = 0.9781117343447195
$ws.Range("D6").Value = 4.417744944668907
$ws.Range("E6").Value = 2.101843225521092
$ws.Range("F6").Value = 2.08518955102823
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.3740519112273883
$ws.Range("C7").Value = 0.9832035832409267
$ws.Range("D7").Value = 4.94152750655546
$ws.Range("E7").Value = 2.222954679375057
$ws.Range("F7").Value = 2.220672418386649
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.3686574517514005
$ws.Range("C8").Value = 0.9897903167014467
$ws.Range("D8").Value = 5.048307810960813
$ws.Range("E8").Value = 2.246843966758887
$ws.Range("F8").Value = 2.246965739777832
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = 0.2563981495475554
$ws.Range("C9").Value = 1.371782596377584
$ws.Range("D9").Value = 8.685257687481315
$ws.Range("E9").Value = 2.947076125158852
$ws.Range("F9").Value = 3.012171393108808
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.4358536319779404
$ws.Range("C10").Value = 1.061643998448061
$ws.Range("D10").Value = 5.46603034801007
$ws.Range("E10").Value = 2.337954308366626
$ws.Range("F10").Value = 2.390760086135826
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.5199476580913862
$ws.Range("C11").Value = 0.536199036767537
$ws.Range("D11").Value = 0.3603403617065257
$ws.Range("E11").Value = 0.600283567746549
$ws.Range("F11").Value = 0.3354004967047021
